# Update "Generate Report for Handback" timestamps for the
# 152c6532-2336-44f0-bf09-f7dd2292f9c2.md row across the Overview,
# zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for 152c6532 row
$overview.Range("G3").Value = "2016-08-26 02:46:16"

# zh-cn sheet: "Correspond Handoff Datetime" / "Correspond Handback DateTime"
$zhcn.Range("H3").Value = "2016-08-26 02:46:11"
$zhcn.Range("K3").Value = "2016-08-26 02:46:30"

# de-de sheet: "Correspond Handoff Datetime" / "Correspond Handback DateTime"
$dede.Range("H3").Value = "2016-08-26 02:46:16"
$dede.Range("K3").Value = "2016-08-26 02:46:37"
